# Apply numeric "want-to-go count" (F) corrections, and one F->"unavailable"
# ticket-price (G) change, across all four worksheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 327
$ws.Range("F3").Value = 1139
$ws.Range("F6").Value = 3392
$ws.Range("F7").Value = 61
$ws.Range("F9").Value = 774
$ws.Range("F10").Value = 595
$ws.Range("F12").Value = 154
$ws.Range("F14").Value = 1814
$ws.Range("F15").Value = 52
$ws.Range("F16").Value = 393
$ws.Range("F17").Value = 53
$ws.Range("F18").Value = 70
$ws.Range("F19").Value = 681
$ws.Range("F20").Value = 451
$ws.Range("F22").Value = 793
$ws.Range("F23").Value = 79964
$ws.Range("F24").Value = 79964
$ws.Range("F27").Value = 33829
$ws.Range("F28").Value = 33829
$ws.Range("F29").Value = 529
$ws.Range("F30").Value = 29
$ws.Range("F31").Value = 26
$ws.Range("F33").Value = 49
$ws.Range("F34").Value = 1000
$ws.Range("F37").Value = 626
$ws.Range("F38").Value = 2490
$ws.Range("F39").Value = 2490
$ws.Range("F40").Value = 1209
$ws.Range("F41").Value = 5497
$ws.Range("F42").Value = 791
$ws.Range("F43").Value = 456
$ws.Range("F47").Value = 411

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G5").Value = "不可售"
$ws.Range("F11").Value = 1977
$ws.Range("F12").Value = 30
$ws.Range("F14").Value = 81
$ws.Range("F16").Value = 10
$ws.Range("F18").Value = 77
$ws.Range("F20").Value = 533
$ws.Range("F21").Value = 533
$ws.Range("F33").Value = 1669
$ws.Range("F35").Value = 10
$ws.Range("F42").Value = 34
$ws.Range("F43").Value = 34
$ws.Range("F47").Value = 192
$ws.Range("F49").Value = 67

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 734
$ws.Range("F5").Value = 579
$ws.Range("F6").Value = 608
$ws.Range("F7").Value = 145

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 734
$ws.Range("F4").Value = 579
$ws.Range("F8").Value = 3392
$ws.Range("F9").Value = 61
$ws.Range("F11").Value = 774
$ws.Range("F12").Value = 608
$ws.Range("F14").Value = 595
$ws.Range("F17").Value = 145
$ws.Range("F18").Value = 1814
$ws.Range("F19").Value = 30
$ws.Range("F20").Value = 52
$ws.Range("F22").Value = 81
$ws.Range("F23").Value = 53
$ws.Range("F24").Value = 70
$ws.Range("F25").Value = 681
$ws.Range("F27").Value = 451
$ws.Range("F29").Value = 79964
$ws.Range("F30").Value = 77
$ws.Range("F31").Value = 33829
$ws.Range("F32").Value = 529
$ws.Range("F33").Value = 26
$ws.Range("F35").Value = 533
$ws.Range("F36").Value = 49
$ws.Range("F38").Value = 1000
$ws.Range("F44").Value = 626
$ws.Range("F45").Value = 2490
$ws.Range("F46").Value = 1209
$ws.Range("F47").Value = 791
$ws.Range("F48").Value = 1669
$ws.Range("F49").Value = 456
$ws.Range("F51").Value = 34
$ws.Range("F52").Value = 34
$ws.Range("F55").Value = 192
